$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 631
$ws.Range("F3").Value = 10732
$ws.Range("F5").Value = 102
$ws.Range("F6").Value = 666
$ws.Range("F8").Value = 12801
$ws.Range("F9").Value = 13225
$ws.Range("F10").Value = 1331
$ws.Range("F11").Value = 1309
$ws.Range("F12").Value = 5571
$ws.Range("F13").Value = 930
$ws.Range("F15").Value = 369
$ws.Range("F17").Value = 1456
$ws.Range("F18").Value = 373
$ws.Range("F19").Value = 2048
$ws.Range("F20").Value = 1071
$ws.Range("F21").Value = 1623
$ws.Range("F23").Value = 21
$ws.Range("F25").Value = 749
$ws.Range("F26").Value = 3074
$ws.Range("F27").Value = 269
$ws.Range("F28").Value = 2123
$ws.Range("F29").Value = 16
$ws.Range("F30").Value = 110
$ws.Range("F31").Value = 1707
$ws.Range("F32").Value = 1019
$ws.Range("F33").Value = 660
$ws.Range("F34").Value = 67
$ws.Range("F36").Value = 3834
$ws.Range("F37").Value = 4491
$ws.Range("F38").Value = 282
$ws.Range("F39").Value = 137
$ws.Range("F42").Value = 3166
$ws.Range("F43").Value = 38
$ws.Range("F45").Value = 311
$ws.Range("F46").Value = 52
$ws.Range("F47").Value = 45
$ws.Range("F48").Value = 4320
$ws.Range("F49").Value = 206

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 90
$ws.Range("F23").Value = 79
$ws.Range("F24").Value = 6

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6767
$ws.Range("F3").Value = 112
$ws.Range("F4").Value = 253

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 632
$ws.Range("F3").Value = 10732
$ws.Range("F5").Value = 112
$ws.Range("F6").Value = 12801
$ws.Range("F7").Value = 13225
$ws.Range("F9").Value = 1331
$ws.Range("F10").Value = 1309
$ws.Range("F11").Value = 5571
$ws.Range("F12").Value = 930
$ws.Range("F13").Value = 369
$ws.Range("F14").Value = 90
$ws.Range("F16").Value = 1456
$ws.Range("F17").Value = 373
$ws.Range("F18").Value = 2048
$ws.Range("F19").Value = 1071
$ws.Range("F20").Value = 1623
$ws.Range("F23").Value = 749
$ws.Range("F24").Value = 3074
$ws.Range("F26").Value = 269
$ws.Range("F27").Value = 2123
$ws.Range("F28").Value = 16
$ws.Range("F29").Value = 110
$ws.Range("F31").Value = 1707
$ws.Range("F33").Value = 1019
$ws.Range("F34").Value = 660
$ws.Range("F35").Value = 67
$ws.Range("F36").Value = 3834
$ws.Range("F37").Value = 4491
$ws.Range("F39").Value = 282
$ws.Range("F40").Value = 137
$ws.Range("F43").Value = 3166
$ws.Range("F45").Value = 311
$ws.Range("F46").Value = 52
$ws.Range("F47").Value = 45
$ws.Range("F48").Value = 4320
$ws.Range("F49").Value = 206
